$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 1.22
$ws.Range("G3").Value = 1.23
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 34
$ws.Range("J3").Value = 6.6
$ws.Range("K3").Value = 7.2
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 3.75
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 1.61
$ws.Range("Q3").Value = 2.42
$ws.Range("R3").Value = 1.15
$ws.Range("S3").Value = 6.6
$ws.Range("T3").Value = 2.44
$ws.Range("U3").Value = 1.51
$ws.Range("V3").Value = 1.03
$ws.Range("W3").Value = 5.3
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 3.8
$ws.Range("AC3").Value = 9.2
$ws.Range("AD3").Value = 46
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 4.5
$ws.Range("AG3").Value = 9.8
$ws.Range("AH3").Value = 48
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 11
$ws.Range("AK3").Value = 26
$ws.Range("AL3").Value = 130
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 24
$ws.Range("AO3").Value = 1000
